# Auto-generated edit script: updates cryptos list (prices/volumes) per commit
# 'Updated cryptos list on Wed Jul 19 16:56:01 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Text)
    $cell = $Sheet.Range($Addr)
    $cell.NumberFormat = '@'
    $cell.Value = $Text
    $cell.NumberFormat = 'General'
}

Set-TextValue $ws 'B21' 'Dai'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D21' '1.001'
Set-TextValue $ws 'E21' '  +0.03%  '

Set-TextValue $ws 'B22' 'WrappedliquidstakedEther2.0'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws 'D22' '2.149.79'
Set-TextValue $ws 'E22' '  +0.60%  '

Set-TextValue $ws 'E23' '  +0.00%  '

Set-TextValue $ws 'D2' '29.911.71'
Set-TextValue $ws 'E2' '  +0.33%  '

Set-TextValue $ws 'D3' '1.903.33'
Set-TextValue $ws 'E3' '  +0.64%  '

Set-TextValue $ws 'E4' '  -0.05%  '

Set-TextValue $ws 'D5' '0.8028'
Set-TextValue $ws 'E5' '  +6.48%  '

Set-TextValue $ws 'D6' '240.60'
Set-TextValue $ws 'E6' '  +0.40%  '

Set-TextValue $ws 'D8' '0.3112'
Set-TextValue $ws 'E8' '  +2.55%  '

Set-TextValue $ws 'D9' '26.29'
Set-TextValue $ws 'E9' '  +3.40%  '

Set-TextValue $ws 'D10' '0.07014'
Set-TextValue $ws 'E10' '  +3.07%  '

Set-TextValue $ws 'D11' '0.07991'
Set-TextValue $ws 'E11' '  +0.54%  '

Set-TextValue $ws 'D12' '1.898.89'
Set-TextValue $ws 'E12' '  +0.41%  '

Set-TextValue $ws 'D13' '0.7389'
Set-TextValue $ws 'E13' '  -0.32%  '

Set-TextValue $ws 'D14' '5.157'
Set-TextValue $ws 'E14' '  -0.14%  '

Set-TextValue $ws 'D15' '92.16'
Set-TextValue $ws 'E15' '  +1.59%  '

Set-TextValue $ws 'D16' '29.905.79'
Set-TextValue $ws 'E16' '  +0.42%  '

Set-TextValue $ws 'D17' '13.93'
Set-TextValue $ws 'E17' '  +0.49%  '

Set-TextValue $ws 'D18' '5.847'
Set-TextValue $ws 'E18' '  -1.60%  '

Set-TextValue $ws 'D19' '244.18'
Set-TextValue $ws 'E19' '  +0.10%  '

Set-TextValue $ws 'D20' '0.000007774'
Set-TextValue $ws 'E20' '  +1.37%  '

Set-TextValue $ws 'D24' '6.885'
Set-TextValue $ws 'E24' '  -0.41%  '

Set-TextValue $ws 'D25' '167.80'
Set-TextValue $ws 'E25' '  +1.46%  '

Set-TextValue $ws 'D26' '9.175'
Set-TextValue $ws 'E26' '  -0.41%  '

Set-TextValue $ws 'D27' '0.1468'
Set-TextValue $ws 'E27' '  +15.32%  '

Set-TextValue $ws 'D28' '18.84'
Set-TextValue $ws 'E28' '  +0.97%  '

Set-TextValue $ws 'E29' '  +2.07%  '

Set-TextValue $ws 'D30' '1.355'
Set-TextValue $ws 'E30' '  -1.89%  '

Set-TextValue $ws 'D31' '1.512'
Set-TextValue $ws 'E31' '  -0.21%  '

Set-TextValue $ws 'D32' '4.278'
Set-TextValue $ws 'E32' '  +0.98%  '

Set-TextValue $ws 'D33' '0.05509'
Set-TextValue $ws 'E33' '  +4.88%  '

Set-TextValue $ws 'D34' '4.046'
Set-TextValue $ws 'E34' '  +1.32%  '

Set-TextValue $ws 'D35' '1.258'
Set-TextValue $ws 'E35' '  +1.24%  '

Set-TextValue $ws 'D36' '0.7293'
Set-TextValue $ws 'E36' '  +1.31%  '

Set-TextValue $ws 'D37' '2.722'
Set-TextValue $ws 'E37' '  +0.17%  '

Set-TextValue $ws 'D38' '0.01914'
Set-TextValue $ws 'E38' '  +0.98%  '

Set-TextValue $ws 'D39' '2.784'
Set-TextValue $ws 'E39' '  +0.42%  '

Set-TextValue $ws 'D40' '0.4396'
Set-TextValue $ws 'E40' '  +0.36%  '

Set-TextValue $ws 'D41' '72.03'
Set-TextValue $ws 'E41' '  +0.49%  '

Set-TextValue $ws 'D42' '5.960'
Set-TextValue $ws 'E42' '  -2.89%  '

Set-TextValue $ws 'D43' '1.001'
Set-TextValue $ws 'E43' '  -0.14%  '

Set-TextValue $ws 'D44' '0.8357'
Set-TextValue $ws 'E44' '  +1.35%  '

Set-TextValue $ws 'E45' '  +0.31%  '

Set-TextValue $ws 'D46' '100.77'
Set-TextValue $ws 'E46' '  +0.29%  '

Set-TextValue $ws 'D47' '7.535'
Set-TextValue $ws 'E47' '  +0.35%  '

Set-TextValue $ws 'D48' '9.686'
Set-TextValue $ws 'E48' '  -0.32%  '

Set-TextValue $ws 'D49' '978.32'
Set-TextValue $ws 'E49' '  +8.74%  '

Set-TextValue $ws 'D50' '2.056.47'
Set-TextValue $ws 'E50' '  +0.33%  '

Set-TextValue $ws 'D51' '36.14'
Set-TextValue $ws 'E51' '  +0.62%  '

